# Branch 05062022: append the newly-received order "2203271DY8MNV3"
# (placed 5/6/2022 4:24:05 AM, Standard Local-Standard Delivery) as five
# line items (qty 1 each) to the bottom of the export sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orderId   = "2203271DY8MNV3"
$qty       = "1"
$courier   = "Standard Local-Standard Delivery"
$timestamp = "5/6/2022 4:24:05 AM"

# SKU (Parent) / CSKU pairs for each line item of this order.
$items = @(
    @("Bath Toys",        "Pink Whale"),
    @("Water Book",       "WB-Happy B-day"),
    @("Water Book",       "WB-Letter Number"),
    @("Water Book",       "WB-Animal"),
    @("Korean Hair Clam", "Coffee Clam")
)

$startRow = 51
$endRow   = $startRow + $items.Length - 1

# Quantity column holds numeric-looking text ("1") throughout the sheet, so
# force the new cells to Text before writing to keep them as strings rather
# than Excel auto-converting them to numbers.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

$row = $startRow
foreach ($item in $items) {
    $ws.Range("A$row").Value = $orderId
    $ws.Range("B$row").Value = $qty
    $ws.Range("C$row").Value = $item[0]
    $ws.Range("D$row").Value = $item[1]
    $ws.Range("E$row").Value = $courier
    $ws.Range("F$row").Value = $timestamp
    $row++
}
